# "extraction effiency stats corrected"
#
# On the "All samples" sheet:
#  - column I ("batch") values for the rows that were "B"/"C" are relabeled
#    ("B"->"C", "C"->"E"); rows that were "A" stay "A"
#  - column J ("extraction") values are switched from letters to roman
#    numerals ("A"->"I", "B"->"II")
#  - a previously-empty "starting_quant" value (column AD) is filled in for
#    rows 42-61
#  - the window scroll/selection is updated to show the bottom of the table

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All samples")
$ws.Activate()

# --- column I ("batch"): relabel the second and third groups of 20 rows ---
$ws.Range("I22:I41").Value = "C"
$ws.Range("I42:I61").Value = "E"

# --- column J ("extraction"): letters -> roman numerals ---
$ws.Range("J2:J11").Value = "I"
$ws.Range("J12:J21").Value = "II"
$ws.Range("J22:J31").Value = "I"
$ws.Range("J32:J41").Value = "II"
$ws.Range("J42:J51").Value = "I"
$ws.Range("J52:J61").Value = "II"

# --- column AD ("starting_quant"): newly-filled values for rows 42-61 ---
$ws.Range("AD42").Value = 154839.733
$ws.Range("AD43").Value = 256626.25
$ws.Range("AD44").Value = 70311.609
$ws.Range("AD45").Value = 71332.6551
$ws.Range("AD46").Value = 31163.203
$ws.Range("AD47").Value = 178117.826
$ws.Range("AD48").Value = 842.412762
$ws.Range("AD49").Value = 51807.2047
$ws.Range("AD50").Value = 145411.677
$ws.Range("AD51").Value = 55663.3961
$ws.Range("AD52").Value = 245859.188
$ws.Range("AD53").Value = 310986.314
$ws.Range("AD54").Value = 58066.6339
$ws.Range("AD55").Value = 73951.0446
$ws.Range("AD56").Value = 24021.5074
$ws.Range("AD57").Value = 257744.585
$ws.Range("AD58").Value = 828.871049
$ws.Range("AD59").Value = 33249.8121
$ws.Range("AD60").Value = 321880.962
$ws.Range("AD61").Value = 65844.3836

# --- window view: re-freeze header row, then move the selection down to
#     the bottom of the table (mirrors the scrolled view in the saved file)
$win = $excel.ActiveWindow
$win.FreezePanes = $false
$ws.Range("A2").Select()
$win.FreezePanes = $true
$ws.Range("J42:J61").Select()
